# Update "想去人数" (want-to-go count) figures in both the "展览" and
# "全部类型" worksheets to reflect freshly scraped totals.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F4").Value = 7280
    $ws.Range("F7").Value = 36
    $ws.Range("F13").Value = 89
    $ws.Range("F14").Value = 660
    $ws.Range("F15").Value = 480
    $ws.Range("F18").Value = 2
    $ws.Range("F20").Value = 71
}
